$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "80391502"
$ws.Range("D16").Value = "IVAN RENE PINILLA ALAYON"
$ws.Range("E16").Value = "1606"
$ws.Range("F16").Value = 27600
$ws.Range("G16").Value = 1300000

$ws.Range("C17").Value = "1047440920"
$ws.Range("D17").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E17").Value = "1609"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 781242

$ws.Range("C18").Value = "1047440920"
$ws.Range("D18").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E18").Value = "1610"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 781242

$ws.Range("C19").Value = "1047440920"
$ws.Range("D19").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E19").Value = "1611"
$ws.Range("F19").Value = 27578
$ws.Range("G19").Value = 781242

$ws.Range("C20").Value = "1047440920"
$ws.Range("D20").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E20").Value = "1612"
$ws.Range("F20").Value = 27578
$ws.Range("G20").Value = 781242

$ws.Range("C21").Value = "1047440920"
$ws.Range("D21").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E21").Value = "1701"
$ws.Range("F21").Value = 27578
$ws.Range("G21").Value = 781242

$ws.Range("C22").Value = "1047440920"
$ws.Range("D22").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E22").Value = "1702"
$ws.Range("F22").Value = 27578
$ws.Range("G22").Value = 781242

$ws.Range("C23").Value = "1047440920"
$ws.Range("D23").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E23").Value = "1703"
$ws.Range("F23").Value = 27578
$ws.Range("G23").Value = 781242

$ws.Range("C24").Value = "1047440920"
$ws.Range("D24").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E24").Value = "1704"
$ws.Range("F24").Value = 27578
$ws.Range("G24").Value = 781242

$ws.Range("C25").Value = "1047440920"
$ws.Range("D25").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E25").Value = "1705"
$ws.Range("F25").Value = 27578
$ws.Range("G25").Value = 781242

$ws.Range("C26").Value = "1047440920"
$ws.Range("D26").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E26").Value = "1706"
$ws.Range("F26").Value = 27578
$ws.Range("G26").Value = 781242

$ws.Range("C27").Value = "1047440920"
$ws.Range("D27").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E27").Value = "1707"
$ws.Range("F27").Value = 27578
$ws.Range("G27").Value = 781242

$ws.Range("C28").Value = "1047440920"
$ws.Range("D28").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E28").Value = "1708"
$ws.Range("F28").Value = 27578
$ws.Range("G28").Value = 781242

$ws.Range("C29").Value = "1047440920"
$ws.Range("D29").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E29").Value = "1709"
$ws.Range("F29").Value = 27578
$ws.Range("G29").Value = 781242

$ws.Range("C30").Value = "1047440920"
$ws.Range("D30").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E30").Value = "1710"
$ws.Range("F30").Value = 27578
$ws.Range("G30").Value = 781242

$ws.Range("C31").Value = "1047440920"
$ws.Range("D31").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E31").Value = "1711"
$ws.Range("F31").Value = 27578
$ws.Range("G31").Value = 781242

$ws.Range("C32").Value = "1047440920"
$ws.Range("D32").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E32").Value = "1712"
$ws.Range("F32").Value = 27578
$ws.Range("G32").Value = 781242

$ws.Range("C33").Value = "1047440920"
$ws.Range("D33").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E33").Value = "1801"
$ws.Range("F33").Value = 27578
$ws.Range("G33").Value = 781242

$ws.Range("C34").Value = "1047440920"
$ws.Range("D34").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E34").Value = "1802"
$ws.Range("F34").Value = 27578
$ws.Range("G34").Value = 781242

$ws.Range("C35").Value = "1047440920"
$ws.Range("D35").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E35").Value = "1803"
$ws.Range("F35").Value = 27578
$ws.Range("G35").Value = 781242

$ws.Range("C36").Value = "1047440920"
$ws.Range("D36").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E36").Value = "1804"
$ws.Range("F36").Value = 27578
$ws.Range("G36").Value = 781242

$ws.Range("C37").Value = "1047440920"
$ws.Range("D37").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E37").Value = "1805"
$ws.Range("F37").Value = 27578
$ws.Range("G37").Value = 781242

$ws.Range("C38").Value = "1047440920"
$ws.Range("D38").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E38").Value = "1806"
$ws.Range("F38").Value = 27578
$ws.Range("G38").Value = 781242

$ws.Range("C39").Value = "1047440920"
$ws.Range("D39").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E39").Value = "1807"
$ws.Range("F39").Value = 27578
$ws.Range("G39").Value = 781242

$ws.Range("C40").Value = "1047440920"
$ws.Range("D40").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E40").Value = "1808"
$ws.Range("F40").Value = 27578
$ws.Range("G40").Value = 781242

$ws.Range("C41").Value = "1047440920"
$ws.Range("D41").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E41").Value = "1809"
$ws.Range("F41").Value = 31249
$ws.Range("G41").Value = 781242

$ws.Range("C42").Value = "1047440920"
$ws.Range("D42").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E42").Value = "1810"
$ws.Range("F42").Value = 31249
$ws.Range("G42").Value = 781242

$ws.Range("C43").Value = "1047440920"
$ws.Range("D43").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E43").Value = "1811"
$ws.Range("F43").Value = 31249
$ws.Range("G43").Value = 781242

$ws.Range("C44").Value = "1047440920"
$ws.Range("D44").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E44").Value = "1812"
$ws.Range("F44").Value = 31249
$ws.Range("G44").Value = 781242

$ws.Range("C45").Value = "1047440920"
$ws.Range("D45").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E45").Value = "1901"
$ws.Range("F45").Value = 31249
$ws.Range("G45").Value = 781242

$ws.Range("C46").Value = "1047440920"
$ws.Range("D46").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E46").Value = "1902"
$ws.Range("F46").Value = 31249
$ws.Range("G46").Value = 781242

$ws.Range("C47").Value = "1047440920"
$ws.Range("D47").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E47").Value = "1903"
$ws.Range("F47").Value = 31249
$ws.Range("G47").Value = 781242

$ws.Range("C48").Value = "1047440920"
$ws.Range("D48").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E48").Value = "1904"
$ws.Range("F48").Value = 31249
$ws.Range("G48").Value = 781242

$ws.Range("C49").Value = "1047440920"
$ws.Range("D49").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E49").Value = "1905"
$ws.Range("F49").Value = 31249
$ws.Range("G49").Value = 781242

$ws.Range("C50").Value = "1047440920"
$ws.Range("D50").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E50").Value = "1906"
$ws.Range("F50").Value = 31249
$ws.Range("G50").Value = 781242

$ws.Range("C51").Value = "1047440920"
$ws.Range("D51").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E51").Value = "1907"
$ws.Range("F51").Value = 31249
$ws.Range("G51").Value = 781242

$ws.Range("C52").Value = "1047440920"
$ws.Range("D52").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E52").Value = "1908"
$ws.Range("F52").Value = 31249
$ws.Range("G52").Value = 781242

$ws.Range("C53").Value = "1047440920"
$ws.Range("D53").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E53").Value = "1909"
$ws.Range("F53").Value = 31249
$ws.Range("G53").Value = 781242

$ws.Range("C54").Value = "1047440920"
$ws.Range("D54").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E54").Value = "1910"
$ws.Range("F54").Value = 31249
$ws.Range("G54").Value = 781242

$ws.Range("C55").Value = "1047440920"
$ws.Range("D55").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E55").Value = "1911"
$ws.Range("F55").Value = 31249
$ws.Range("G55").Value = 781242

$ws.Range("C56").Value = "1047440920"
$ws.Range("D56").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E56").Value = "1912"
$ws.Range("F56").Value = 31249
$ws.Range("G56").Value = 781242

$ws.Range("C57").Value = "1047440920"
$ws.Range("D57").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E57").Value = "2001"
$ws.Range("F57").Value = 31249
$ws.Range("G57").Value = 781242

$ws.Range("C58").Value = "1047440920"
$ws.Range("D58").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E58").Value = "2002"
$ws.Range("F58").Value = 31249
$ws.Range("G58").Value = 781242

$ws.Range("C59").Value = "1047440920"
$ws.Range("D59").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E59").Value = "2003"
$ws.Range("F59").Value = 31249
$ws.Range("G59").Value = 781242

$ws.Range("C60").Value = "1047440920"
$ws.Range("D60").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E60").Value = "2004"
$ws.Range("F60").Value = 31249
$ws.Range("G60").Value = 781242

$ws.Range("C61").Value = "1047440920"
$ws.Range("D61").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E61").Value = "2005"
$ws.Range("F61").Value = 31249
$ws.Range("G61").Value = 781242

$ws.Range("C62").Value = "1047440920"
$ws.Range("D62").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E62").Value = "2006"
$ws.Range("F62").Value = 31249
$ws.Range("G62").Value = 781242

$ws.Range("C63").Value = "1047440920"
$ws.Range("D63").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E63").Value = "2007"
$ws.Range("F63").Value = 31249
$ws.Range("G63").Value = 781242

$ws.Range("C64").Value = "1047440920"
$ws.Range("D64").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E64").Value = "2008"
$ws.Range("F64").Value = 31249
$ws.Range("G64").Value = 781242

$ws.Range("C65").Value = "1047440920"
$ws.Range("D65").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E65").Value = "2009"
$ws.Range("F65").Value = 31249
$ws.Range("G65").Value = 781242

$ws.Range("C66").Value = "1047440920"
$ws.Range("D66").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E66").Value = "2010"
$ws.Range("F66").Value = 31249
$ws.Range("G66").Value = 781242

$ws.Range("C67").Value = "1047440920"
$ws.Range("D67").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E67").Value = "2011"
$ws.Range("F67").Value = 31249
$ws.Range("G67").Value = 781242

$ws.Range("C68").Value = "1047440920"
$ws.Range("D68").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E68").Value = "2012"
$ws.Range("F68").Value = 31249
$ws.Range("G68").Value = 781242

$ws.Range("C69").Value = "1047440920"
$ws.Range("D69").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E69").Value = "2101"
$ws.Range("F69").Value = 31249
$ws.Range("G69").Value = 781242

$ws.Range("C70").Value = "1047440920"
$ws.Range("D70").Value = "BAYRON CASTRO CARDENAS"
$ws.Range("E70").Value = "2102"
$ws.Range("F70").Value = 22916
$ws.Range("G70").Value = 781242
